$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.065987540171553
$ws.Cells.Item(2, 4).Value = 1.06626980592657
$ws.Cells.Item(2, 5).Value = 1.069694012749479
$ws.Cells.Item(2, 6).Value = 1.07548978111005
$ws.Cells.Item(2, 9).Value = 1.051618085871973
$ws.Cells.Item(2, 10).Value = 1.070939845295699
$ws.Cells.Item(2, 11).Value = 1.068980807371998
$ws.Cells.Item(2, 12).Value = 1.072395836392625
$ws.Cells.Item(2, 13).Value = 1.078176214448922
$ws.Cells.Item(2, 14).Value = 1.072460702714536

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.067374068262262
$ws.Cells.Item(3, 4).Value = 1.067368066606238
$ws.Cells.Item(3, 5).Value = 1.071032202015362
$ws.Cells.Item(3, 6).Value = 1.07679624286016
$ws.Cells.Item(3, 9).Value = 1.052063197447674
$ws.Cells.Item(3, 10).Value = 1.071980122685763
$ws.Cells.Item(3, 11).Value = 1.069893922464409
$ws.Cells.Item(3, 12).Value = 1.073548951770737
$ws.Cells.Item(3, 13).Value = 1.079298802755824
$ws.Cells.Item(3, 14).Value = 1.073502457417816

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.068270349286634
$ws.Cells.Item(4, 4).Value = 1.068077858426915
$ws.Cells.Item(4, 5).Value = 1.071897537051234
$ws.Cells.Item(4, 6).Value = 1.077641118955275
$ws.Cells.Item(4, 9).Value = 1.05234949458882
$ws.Cells.Item(4, 10).Value = 1.072651872903893
$ws.Cells.Item(4, 11).Value = 1.070483319425082
$ws.Cells.Item(4, 12).Value = 1.074293986302073
$ws.Cells.Item(4, 13).Value = 1.08002414636963
$ws.Cells.Item(4, 14).Value = 1.074175161598305

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.068646936317945
$ws.Cells.Item(5, 4).Value = 1.068376052968287
$ws.Cells.Item(5, 5).Value = 1.07226119321984
$ws.Cells.Item(5, 6).Value = 1.077996190816976
$ws.Cells.Item(5, 9).Value = 1.052469443809223
$ws.Cells.Item(5, 10).Value = 1.072933950045571
$ws.Cells.Item(5, 11).Value = 1.070730757796366
$ws.Cells.Item(5, 12).Value = 1.074606937206721
$ws.Cells.Item(5, 13).Value = 1.080328833808568
$ws.Cells.Item(5, 14).Value = 1.074457639321881

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.068710154792122
$ws.Cells.Item(6, 4).Value = 1.06842610939057
$ws.Cells.Item(6, 5).Value = 1.072322245145835
$ws.Cells.Item(6, 6).Value = 1.078055802365557
$ws.Cells.Item(6, 9).Value = 1.052489559810017
$ws.Cells.Item(6, 10).Value = 1.072981292983704
$ws.Cells.Item(6, 11).Value = 1.070772283693411
$ws.Cells.Item(6, 12).Value = 1.074659467806836
$ws.Cells.Item(6, 13).Value = 1.080379977784216
$ws.Cells.Item(6, 14).Value = 1.074505049492416

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.068275382078193
$ws.Cells.Item(7, 4).Value = 1.068081843706666
$ws.Cells.Item(7, 5).Value = 1.071902396751169
$ws.Cells.Item(7, 6).Value = 1.077645863886351
$ws.Cells.Item(7, 9).Value = 1.052351098965558
$ws.Cells.Item(7, 10).Value = 1.072655643311706
$ws.Cells.Item(7, 11).Value = 1.070486627057213
$ws.Cells.Item(7, 12).Value = 1.07429816899074
$ws.Cells.Item(7, 13).Value = 1.080028218584462
$ws.Cells.Item(7, 14).Value = 1.07417893736053

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.066456312140949
$ws.Cells.Item(8, 4).Value = 1.066641147405873
$ws.Cells.Item(8, 5).Value = 1.070146379089376
$ws.Cells.Item(8, 6).Value = 1.07593141010705
$ws.Cells.Item(8, 9).Value = 1.0517688704938
$ws.Cells.Item(8, 10).Value = 1.071291699477972
$ws.Cells.Item(8, 11).Value = 1.069289700793611
$ws.Cells.Item(8, 12).Value = 1.07278576882935
$ws.Cells.Item(8, 13).Value = 1.07855581708401
$ws.Cells.Item(8, 14).Value = 1.072813056570106

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.063243803856691
$ws.Cells.Item(9, 4).Value = 1.064095760972361
$ws.Cells.Item(9, 5).Value = 1.0670475572455
$ws.Cells.Item(9, 6).Value = 1.072906387111666
$ws.Cells.Item(9, 9).Value = 1.050729664568405
$ws.Cells.Item(9, 10).Value = 1.068877541495544
$ws.Cells.Item(9, 11).Value = 1.067169330713042
$ws.Cells.Item(9, 12).Value = 1.070112072660167
$ws.Cells.Item(9, 13).Value = 1.075953092094771
$ws.Cells.Item(9, 14).Value = 1.070395470206434

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.06109706959807
$ws.Cells.Item(10, 4).Value = 1.062394140646258
$ws.Cells.Item(10, 5).Value = 1.064978396004483
$ws.Cells.Item(10, 6).Value = 1.070886821236318
$ws.Cells.Item(10, 9).Value = 1.050027851800133
$ws.Cells.Item(10, 10).Value = 1.067260680819452
$ws.Cells.Item(10, 11).Value = 1.065748012889875
$ws.Cells.Item(10, 12).Value = 1.068323555866997
$ws.Cells.Item(10, 13).Value = 1.074212235233164
$ws.Cells.Item(10, 14).Value = 1.068776313402727

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.060166240789177
$ws.Cells.Item(11, 4).Value = 1.06165616247825
$ws.Cells.Item(11, 5).Value = 1.064081590261484
$ws.Cells.Item(11, 6).Value = 1.070011589421878
$ws.Cells.Item(11, 9).Value = 1.049721798899293
$ws.Cells.Item(11, 10).Value = 1.066558752775639
$ws.Cells.Item(11, 11).Value = 1.065130689465994
$ws.Cells.Item(11, 12).Value = 1.067547619982952
$ws.Cells.Item(11, 13).Value = 1.073457022497018
$ws.Cells.Item(11, 14).Value = 1.068073388540579

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.059820291427459
$ws.Cells.Item(12, 4).Value = 1.061381865572649
$ws.Cells.Item(12, 5).Value = 1.06374834459341
$ws.Cells.Item(12, 6).Value = 1.069686372512769
$ws.Cells.Item(12, 9).Value = 1.049607790121942
$ws.Cells.Item(12, 10).Value = 1.066297748188992
$ws.Cells.Item(12, 11).Value = 1.064901101407074
$ws.Cells.Item(12, 12).Value = 1.067259173410919
$ws.Cells.Item(12, 13).Value = 1.073176286688106
$ws.Cells.Item(12, 14).Value = 1.067812013297482

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.059894507800647
$ws.Cells.Item(13, 4).Value = 1.061440711350643
$ws.Cells.Item(13, 5).Value = 1.063819832979633
$ws.Cells.Item(13, 6).Value = 1.069756138006854
$ws.Cells.Item(13, 9).Value = 1.049632260227403
$ws.Cells.Item(13, 10).Value = 1.066353747174938
$ws.Cells.Item(13, 11).Value = 1.064950361856186
$ws.Cells.Item(13, 12).Value = 1.067321056636905
$ws.Cells.Item(13, 13).Value = 1.073236515323718
$ws.Cells.Item(13, 14).Value = 1.067868091808411

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.060137648582869
$ws.Cells.Item(14, 4).Value = 1.061633492683613
$ws.Cells.Item(14, 5).Value = 1.064054046794654
$ws.Cells.Item(14, 6).Value = 1.069984709301548
$ws.Cells.Item(14, 9).Value = 1.049712381583098
$ws.Cells.Item(14, 10).Value = 1.066537183736131
$ws.Cells.Item(14, 11).Value = 1.065111717514473
$ws.Cells.Item(14, 12).Value = 1.067523781606642
$ws.Cells.Item(14, 13).Value = 1.07343382122231
$ws.Cells.Item(14, 14).Value = 1.068051788870561

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.060287429214104
$ws.Cells.Item(15, 4).Value = 1.061752247795723
$ws.Cells.Item(15, 5).Value = 1.064198336000774
$ws.Cells.Item(15, 6).Value = 1.070125523996182
$ws.Cells.Item(15, 9).Value = 1.049761703587781
$ws.Cells.Item(15, 10).Value = 1.066650168186687
$ws.Cells.Item(15, 11).Value = 1.065211095965924
$ws.Cells.Item(15, 12).Value = 1.06764865660973
$ws.Cells.Item(15, 13).Value = 1.073555359124232
$ws.Cells.Item(15, 14).Value = 1.068164933771996

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.06115881818941
$ws.Cells.Item(16, 4).Value = 1.06244309294412
$ws.Cells.Item(16, 5).Value = 1.065037895805182
$ws.Cells.Item(16, 6).Value = 1.070944891362013
$ws.Cells.Item(16, 9).Value = 1.050048117766351
$ws.Cells.Item(16, 10).Value = 1.067307226828071
$ws.Cells.Item(16, 11).Value = 1.065788942616353
$ws.Cells.Item(16, 12).Value = 1.068375020255788
$ws.Cells.Item(16, 13).Value = 1.074262326136714
$ws.Cells.Item(16, 14).Value = 1.068822925512017

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.06170507090889
$ws.Cells.Item(17, 4).Value = 1.062876127129971
$ws.Cells.Item(17, 5).Value = 1.065564299599751
$ws.Cells.Item(17, 6).Value = 1.071458655706513
$ws.Cells.Item(17, 9).Value = 1.050227197156874
$ws.Cells.Item(17, 10).Value = 1.067718892930027
$ws.Cells.Item(17, 11).Value = 1.066150903708438
$ws.Cells.Item(17, 12).Value = 1.068830245113175
$ws.Cells.Item(17, 13).Value = 1.074705407277966
$ws.Cells.Item(17, 14).Value = 1.069235176227057

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.062023567908183
$ws.Cells.Item(18, 4).Value = 1.063128596470008
$ws.Cells.Item(18, 5).Value = 1.06587126081701
$ws.Cells.Item(18, 6).Value = 1.071758254053095
$ws.Cells.Item(18, 9).Value = 1.050331442546469
$ws.Cells.Item(18, 10).Value = 1.067958835810878
$ws.Cells.Item(18, 11).Value = 1.066361848017199
$ws.Cells.Item(18, 12).Value = 1.069095626013532
$ws.Cells.Item(18, 13).Value = 1.074963713237568
$ws.Cells.Item(18, 14).Value = 1.06947545985432

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.062132146476618
$ws.Cells.Item(19, 4).Value = 1.063214663080425
$ws.Cells.Item(19, 5).Value = 1.065975913038112
$ws.Cells.Item(19, 6).Value = 1.07186039732609
$ws.Cells.Item(19, 9).Value = 1.050366952212771
$ws.Cells.Item(19, 10).Value = 1.068040620580924
$ws.Cells.Item(19, 11).Value = 1.066433743905312
$ws.Cells.Item(19, 12).Value = 1.069186089747302
$ws.Cells.Item(19, 13).Value = 1.075051765994406
$ws.Cells.Item(19, 14).Value = 1.06955736076812

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.061646475926066
$ws.Cells.Item(20, 4).Value = 1.062829678292002
$ws.Cells.Item(20, 5).Value = 1.065507829892689
$ws.Cells.Item(20, 6).Value = 1.071403541106572
$ws.Cells.Item(20, 9).Value = 1.050208005224157
$ws.Cells.Item(20, 10).Value = 1.067674743185967
$ws.Cells.Item(20, 11).Value = 1.066112087507174
$ws.Cells.Item(20, 12).Value = 1.068781418729556
$ws.Cells.Item(20, 13).Value = 1.074657882940926
$ws.Cells.Item(20, 14).Value = 1.069190963785296

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.0600660551993
$ws.Cells.Item(21, 4).Value = 1.061576728357305
$ws.Cells.Item(21, 5).Value = 1.063985080350985
$ws.Cells.Item(21, 6).Value = 1.069917404015553
$ws.Cells.Item(21, 9).Value = 1.049688796887267
$ws.Cells.Item(21, 10).Value = 1.066483173929376
$ws.Cells.Item(21, 11).Value = 1.065064210237552
$ws.Cells.Item(21, 12).Value = 1.067464090526736
$ws.Cells.Item(21, 13).Value = 1.073375725554998
$ws.Cells.Item(21, 14).Value = 1.067997702363685

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.059071232847434
$ws.Cells.Item(22, 4).Value = 1.060787911685479
$ws.Cells.Item(22, 5).Value = 1.063026900531945
$ws.Cells.Item(22, 6).Value = 1.06898233227643
$ws.Cells.Item(22, 9).Value = 1.049360456440681
$ws.Cells.Item(22, 10).Value = 1.065732380760083
$ws.Cells.Item(22, 11).Value = 1.064403707880047
$ws.Cells.Item(22, 12).Value = 1.066634504241594
$ws.Cells.Item(22, 13).Value = 1.072568329460162
$ws.Cells.Item(22, 14).Value = 1.067245842981975

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.05959871792722
$ws.Cells.Item(23, 4).Value = 1.06120617795378
$ws.Cells.Item(23, 5).Value = 1.063534924168699
$ws.Cells.Item(23, 6).Value = 1.069478097372721
$ws.Cells.Item(23, 9).Value = 1.049534696110117
$ws.Cells.Item(23, 10).Value = 1.0661305441294
$ws.Cells.Item(23, 11).Value = 1.064754011293493
$ws.Cells.Item(23, 12).Value = 1.06707441122612
$ws.Cells.Item(23, 13).Value = 1.072996465604966
$ws.Cells.Item(23, 14).Value = 1.067644571788946

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.0616729528562
$ws.Cells.Item(24, 4).Value = 1.06285066686893
$ws.Cells.Item(24, 5).Value = 1.065533346373249
$ws.Cells.Item(24, 6).Value = 1.07142844524322
$ws.Cells.Item(24, 9).Value = 1.050216677876839
$ws.Cells.Item(24, 10).Value = 1.067694693095354
$ws.Cells.Item(24, 11).Value = 1.066129627438222
$ws.Cells.Item(24, 12).Value = 1.068803481715459
$ws.Cells.Item(24, 13).Value = 1.074679357562464
$ws.Cells.Item(24, 14).Value = 1.069210942025842

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.064075183995139
$ws.Cells.Item(25, 4).Value = 1.064754617766779
$ws.Cells.Item(25, 5).Value = 1.067849236082261
$ws.Cells.Item(25, 6).Value = 1.073688918897511
$ws.Cells.Item(25, 9).Value = 1.050999904443262
$ws.Cells.Item(25, 10).Value = 1.069502951326533
$ws.Cells.Item(25, 11).Value = 1.06771884739313
$ws.Cells.Item(25, 12).Value = 1.07080433678263
$ws.Cells.Item(25, 13).Value = 1.076626948794407
$ws.Cells.Item(25, 14).Value = 1.07102176819112

